$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "SVMXC__Service_Order__c WO = new SVMXC__Service_Order__c ( SVMXC__Company__c = '001q000000kxZfw', SVMXC__Order_Status__c = 'Open' ,SVMXC__Country__c = 'United States', SVMXC__Sub_Status__c = 'Resolved', SVMXC__Priority__c = 'High',SVMXC__Actual_Initial_Response__c = System.Today(),SVMXC__Actual_Onsite_Response__c = System.Today() - 1 );insert WO;"
$ws.Range("B2").Value = "Select Name , Id,FORMAT(SVMXC__Actual_Initial_Response__c), FORMAT(SVMXC__Actual_Onsite_Response__c) from SVMXC__Service_Order__c where Createdby.Id = '005q0000003GGfP' Order by CreatedDate DESC Limit 1"

$ws.Range("B11").Select()
